# Applies the cryptos-list data refresh described in the commit message.
# Column D values that look like plain numbers (e.g. "1.003") are written
# with a leading apostrophe so Excel stores them as text (matching the
# original inline-string cell content) instead of coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.375.36"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.831.38"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'252.34"
$ws.Range("E5").Value = "  -2.93%  "

$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "'0.5232"
$ws.Range("E7").Value = "  -0.55%  "

$ws.Range("D8").Value = "'0.2760"
$ws.Range("E8").Value = "  -13.78%  "

$ws.Range("D9").Value = "'0.06813"
$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("D10").Value = "1.853.07"
$ws.Range("E10").Value = "  +0.82%  "

$ws.Range("D11").Value = "'16.45"
$ws.Range("E11").Value = "  -12.23%  "

$ws.Range("D12").Value = "'0.07102"
$ws.Range("E12").Value = "  -8.19%  "

$ws.Range("E13").Value = "  -12.70%  "

$ws.Range("D14").Value = "'85.74"
$ws.Range("E14").Value = "  -2.10%  "

$ws.Range("D15").Value = "'4.837"
$ws.Range("E15").Value = "  -3.36%  "

$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "'13.14"
$ws.Range("E18").Value = "  -4.95%  "

$ws.Range("D19").Value = "'0.000007316"
$ws.Range("E19").Value = "  -7.75%  "

$ws.Range("D20").Value = "26.394.75"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").Value = "2.085.04"
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").Value = "'4.484"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("D23").Value = "'5.793"
$ws.Range("E23").Value = "  -2.97%  "

$ws.Range("D24").Value = "'8.934"
$ws.Range("E24").Value = "  -4.63%  "

$ws.Range("D25").Value = "'142.01"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").Value = "'1.676"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'2.021"
$ws.Range("E27").Value = "  -5.74%  "

$ws.Range("D28").Value = "'16.47"
$ws.Range("E28").Value = "  -2.63%  "

$ws.Range("D29").Value = "'108.79"
$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("D30").Value = "'4.038"
$ws.Range("E30").Value = "  -2.40%  "

$ws.Range("D32").Value = "'3.831"
$ws.Range("E32").Value = "  -5.53%  "

$ws.Range("D33").Value = "'0.04688"
$ws.Range("E33").Value = "  -3.48%  "

$ws.Range("D34").Value = "'2.878"
$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6995"
$ws.Range("E35").Value = "  -3.99%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.096"
$ws.Range("E36").Value = "  -3.02%  "

$ws.Range("D37").Value = "'3.053"
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("D38").Value = "'2.161"
$ws.Range("E38").Value = "  -3.40%  "

$ws.Range("D39").Value = "'0.01631"
$ws.Range("E39").Value = "  -7.37%  "

$ws.Range("D40").Value = "'0.4417"
$ws.Range("E40").Value = "  -7.39%  "

$ws.Range("D41").Value = "'0.8602"
$ws.Range("E41").Value = "  -3.55%  "

$ws.Range("D42").Value = "'104.82"
$ws.Range("E42").Value = "  -4.33%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "'5.714"
$ws.Range("E44").Value = "  -3.18%  "

$ws.Range("D45").Value = "'6.960"
$ws.Range("E45").Value = "  -8.89%  "

$ws.Range("D46").Value = "'8.574"
$ws.Range("E46").Value = "  -3.95%  "

$ws.Range("D47").Value = "'0.05570"
$ws.Range("E47").Value = "  -4.79%  "

$ws.Range("D48").Value = "'58.44"
$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'33.18"
$ws.Range("E49").Value = "  -4.74%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'0.8570"
$ws.Range("E50").Value = "  -4.18%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1167"
$ws.Range("E51").Value = "  -5.10%  "
